$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: "?????" (highlighted red) -> "April 31,2024" (no highlight)
# -----------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("?????", $true, $false, $false, $false, $false, $true, 1, $false, "April 31,2024", 2)

# Re-find the freshly inserted text so we can clear its highlight.
$rng1b = $d.Content
$rng1b.Find.Execute("April 31,2024", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng1b.HighlightColorIndex = 0

# -----------------------------------------------------------------
# Change 2: "28th" + " " + "Feb" + " 202" + "3" (five runs) ->
#           a single run reading "April 31,2024"
# -----------------------------------------------------------------
# Locate the run immediately preceding "28th Feb 2023" ("... prior to ")
# so we can briefly perturb its formatting; this stops the editor from
# folding that unrelated run into the newly written text (both runs
# otherwise share identical run formatting). Match the run's full text
# (including its leading/trailing spaces) so the whole run is selected.
$before2 = $d.Content
$before2.Find.Execute(" prior to ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$priorRun = $d.Range($before2.Start, $before2.End)

$origSize = $priorRun.Font.Size
$tempSize = $origSize - 1
$priorRun.Font.Size = $tempSize

# Now replace the five-run date phrase with the corrected single string.
$rng2 = $d.Content
$rng2.Find.Execute("28th Feb 2023", $true, $false, $false, $false, $false, $true, 1, $false, "April 31,2024", 2)

# Restore the preceding run's formatting exactly as it was.
$priorAgain = $d.Content
$priorAgain.Find.Execute(" prior to ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$priorRunFixed = $d.Range($priorAgain.Start, $priorAgain.End)
$priorRunFixed.Font.Size = $origSize
